$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,2).Value = "''Akkermansia_muciniphila_ATCC_BAA_835.mat'"
$ws.Cells.Item(2,2).Style = "Normal"
$ws.Cells.Item(2,3).Value = 0

$ws.Cells.Item(3,2).Value = "''Alistipes_finegoldii_DSM_17242.mat'"
$ws.Cells.Item(3,2).Style = "Normal"
$ws.Cells.Item(3,3).Value = 0

$ws.Cells.Item(4,2).Value = "''Alistipes_indistinctus_YIT_12060.mat'"
$ws.Cells.Item(4,2).Style = "Normal"
$ws.Cells.Item(4,3).Value = 0

$ws.Cells.Item(5,2).Value = "''Alistipes_putredinis_DSM_17216.mat'"
$ws.Cells.Item(5,2).Style = "Normal"
$ws.Cells.Item(5,3).Value = 0

$ws.Cells.Item(6,2).Value = "''Alistipes_shahii_WAL_8301.mat'"
$ws.Cells.Item(6,2).Style = "Normal"
$ws.Cells.Item(6,3).Value = 0

$ws.Cells.Item(7,2).Value = "''Bacteroides_cellulosilyticus_DSM_14838.mat'"
$ws.Cells.Item(7,2).Style = "Normal"
$ws.Cells.Item(7,3).Value = 0

$ws.Cells.Item(8,2).Value = "''Bacteroides_coprophilus_DSM_18228.mat'"
$ws.Cells.Item(8,2).Style = "Normal"
$ws.Cells.Item(8,3).Value = 0

$ws.Cells.Item(9,2).Value = "''Bacteroides_fragilis_3_1_12.mat'"
$ws.Cells.Item(9,2).Style = "Normal"
$ws.Cells.Item(9,3).Value = 0

$ws.Cells.Item(10,2).Value = "''Bacteroides_oleiciplenus_YIT_12058.mat'"
$ws.Cells.Item(10,2).Style = "Normal"
$ws.Cells.Item(10,3).Value = 0

$ws.Cells.Item(11,2).Value = "''Bacteroides_ovatus_ATCC_8483.mat'"
$ws.Cells.Item(11,2).Style = "Normal"
$ws.Cells.Item(11,3).Value = 0

$ws.Cells.Item(12,2).Value = "''Bacteroides_plebeius_M12_DSM_17135.mat'"
$ws.Cells.Item(12,2).Style = "Normal"
$ws.Cells.Item(12,3).Value = 0

$ws.Cells.Item(13,2).Value = "''Bacteroides_salyersiae_WAL_10018.mat'"
$ws.Cells.Item(13,2).Style = "Normal"
$ws.Cells.Item(13,3).Value = 0

$ws.Cells.Item(14,2).Value = "''Bacteroides_thetaiotaomicron_VPI_5482.mat'"
$ws.Cells.Item(14,2).Style = "Normal"
$ws.Cells.Item(14,3).Value = 0

$ws.Cells.Item(15,2).Value = "''Bacteroides_uniformis_ATCC_8492.mat'"
$ws.Cells.Item(15,2).Style = "Normal"
$ws.Cells.Item(15,3).Value = 0

$ws.Cells.Item(16,2).Value = "''Bacteroides_vulgatus_ATCC_8482.mat'"
$ws.Cells.Item(16,2).Style = "Normal"
$ws.Cells.Item(16,3).Value = 0

$ws.Cells.Item(17,2).Value = "''Barnesiella_intestinihominis_YIT_11860.mat'"
$ws.Cells.Item(17,2).Style = "Normal"
$ws.Cells.Item(17,3).Value = 0

$ws.Cells.Item(18,2).Value = "''Bifidobacterium_animalis_lactis_AD011.mat'"
$ws.Cells.Item(18,2).Style = "Normal"
$ws.Cells.Item(18,3).Value = 0

$ws.Cells.Item(19,2).Value = "''Bilophila_wadsworthia_3_1_6.mat'"
$ws.Cells.Item(19,2).Style = "Normal"
$ws.Cells.Item(19,3).Value = 0

$ws.Cells.Item(20,2).Value = "''Escherichia_coli_O157_H7_str_Sakai.mat'"
$ws.Cells.Item(20,2).Style = "Normal"
$ws.Cells.Item(20,3).Value = 0

$ws.Cells.Item(21,2).Value = "''Eubacterium_ramulus_ATCC_29099.mat'"
$ws.Cells.Item(21,2).Style = "Normal"
$ws.Cells.Item(21,3).Value = 0

$ws.Cells.Item(22,2).Value = "''Flavonifractor_plautii_ATCC_29863.mat'"
$ws.Cells.Item(22,2).Style = "Normal"
$ws.Cells.Item(22,3).Value = 0.004

$ws.Cells.Item(23,2).Value = "''Marvinbryantia_formatexigens_I_52_DSM_14469.mat'"
$ws.Cells.Item(23,2).Style = "Normal"
$ws.Cells.Item(23,3).Value = 0

$ws.Cells.Item(24,2).Value = "''Odoribacter_splanchnicus_1651_6_DSM_20712.mat'"
$ws.Cells.Item(24,2).Style = "Normal"
$ws.Cells.Item(24,3).Value = 0

$ws.Cells.Item(25,2).Value = "''Parabacteroides_distasonis_ATCC_8503.mat'"
$ws.Cells.Item(25,2).Style = "Normal"
$ws.Cells.Item(25,3).Value = 0

$ws.Cells.Item(26,2).Value = "''Parabacteroides_johnsonii_DSM_18315.mat'"
$ws.Cells.Item(26,2).Style = "Normal"
$ws.Cells.Item(26,3).Value = 0

$ws.Cells.Item(27,2).Value = "''Paraprevotella_xylaniphila_YIT_11841.mat'"
$ws.Cells.Item(27,2).Style = "Normal"
$ws.Cells.Item(27,3).Value = 0

$ws.Cells.Item(28,2).Value = "''Parasutterella_excrementihominis_YIT_11859.mat'"
$ws.Cells.Item(28,2).Style = "Normal"
$ws.Cells.Item(28,3).Value = 0.032

$ws.Cells.Item(29,2).Value = "''Phascolarctobacterium_succinatutens_YIT_12067.mat'"
$ws.Cells.Item(29,2).Style = "Normal"
$ws.Cells.Item(29,3).Value = 0.075

$ws.Cells.Item(30,2).Value = "''Prevotella_copri_CB7_DSM_18205.mat'"
$ws.Cells.Item(30,2).Style = "Normal"
$ws.Cells.Item(30,3).Value = 0

$ws.Cells.Item(31,2).Value = "''Prevotella_stercorea_DSM_18206.mat'"
$ws.Cells.Item(31,2).Style = "Normal"
$ws.Cells.Item(31,3).Value = 0.023

$ws.Cells.Item(32,2).Value = "''Roseburia_inulinivorans_DSM_16841.mat'"
$ws.Cells.Item(32,2).Style = "Normal"
$ws.Cells.Item(32,3).Value = 0.133

$ws.Cells.Item(33,2).Value = "''Sutterella_wadsworthensis_3_1_45B.mat'"
$ws.Cells.Item(33,2).Style = "Normal"
$ws.Cells.Item(33,3).Value = 0.732
